$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")

# Update Version value
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for "Jurisdiction"
$meta.Rows.Item(11).Insert()
$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# Copy formatting from the row above (Contact row) into the new row so the
# new cells keep the same body style as the rest of the table
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122) | Out-Null
